# Auto-generated script applying scheduled price-refresh updates
# to the Leve profit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 324.5
$ws.Range("I11").Value = 324.5
$ws.Range("K11").Value = 324.5
$ws.Range("M11").Value = -184.5
$ws.Range("H17").Value = 2179.1765
$ws.Range("J17").Value = 2179.1765
$ws.Range("L17").Value = 6537.529500000001
$ws.Range("N17").Value = -6873.529500000001
$ws.Range("H46").Value = 3908.5
$ws.Range("I46").Value = 17.0
$ws.Range("J46").Value = 7800.0
$ws.Range("K46").Value = 51.0
$ws.Range("L46").Value = 23400.0
$ws.Range("M46").Value = 68.0
$ws.Range("N46").Value = -23638.0
$ws.Range("H60").Value = 3908.5
$ws.Range("I60").Value = 17.0
$ws.Range("J60").Value = 7800.0
$ws.Range("K60").Value = 51.0
$ws.Range("L60").Value = 23400.0
$ws.Range("M60").Value = 433.0
$ws.Range("N60").Value = -24368.0
$ws.Range("H98").Value = 1041.9412
$ws.Range("I98").Value = 876.3571
$ws.Range("K98").Value = 876.3571
$ws.Range("M98").Value = 621.6429
$ws.Range("H116").Value = 12257.714
$ws.Range("I116").Value = 14099.0
$ws.Range("J116").Value = 11521.2
$ws.Range("K116").Value = 14099.0
$ws.Range("L116").Value = 11521.2
$ws.Range("M116").Value = -10657.0
$ws.Range("N116").Value = -18405.2
$ws.Range("H122").Value = 1041.9412
$ws.Range("I122").Value = 876.3571
$ws.Range("K122").Value = 2629.0713
$ws.Range("M122").Value = -179.0712999999996
$ws.Range("H127").Value = 2727.2163
$ws.Range("I127").Value = 981.8
$ws.Range("K127").Value = 2945.4
$ws.Range("M127").Value = 2014.6
$ws.Range("H134").Value = 125000.0
$ws.Range("J134").Value = 125000.0
$ws.Range("L134").Value = 125000.0
$ws.Range("N134").Value = -135140.0
$ws.Range("H136").Value = 101601.1
$ws.Range("J136").Value = 101601.1
$ws.Range("L136").Value = 101601.1
$ws.Range("N136").Value = -111801.1
$ws.Range("H137").Value = 2600.9333
$ws.Range("I137").Value = 2600.9333
$ws.Range("J137").Value = 0.0
$ws.Range("K137").Value = 7802.7999
$ws.Range("L137").Value = 0.0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25420.816
$ws.Range("I32").Value = 27344.654
$ws.Range("K32").Value = 27344.654
$ws.Range("M32").Value = -27057.654
$ws.Range("H102").Value = 2636.1667
$ws.Range("I102").Value = 3332.6667
$ws.Range("J102").Value = 1939.6666
$ws.Range("K102").Value = 3332.6667
$ws.Range("L102").Value = 1939.6666
$ws.Range("M102").Value = -1710.6667
$ws.Range("N102").Value = -5183.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3177.261
$ws.Range("I105").Value = 3172.5264
$ws.Range("K105").Value = 3172.5264
$ws.Range("M105").Value = -1425.5264
$ws.Range("H134").Value = 3633.5557
$ws.Range("I134").Value = 3092.6667
$ws.Range("K134").Value = 9278.000100000001
$ws.Range("M134").Value = -6743.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0.0
$ws.Range("I57").Value = 0.0
$ws.Range("K57").Value = 0.0
$ws.Range("M57").ClearContents()
$ws.Range("H104").Value = 25000.0
$ws.Range("I104").Value = 0.0
$ws.Range("K104").Value = 0.0
$ws.Range("M104").ClearContents()
$ws.Range("H107").Value = 1309.0
$ws.Range("I107").Value = 670.0
$ws.Range("J107").Value = 2906.5
$ws.Range("K107").Value = 670.0
$ws.Range("L107").Value = 2906.5
$ws.Range("M107").Value = 1250.0
$ws.Range("N107").Value = -6746.5
$ws.Range("H132").Value = 43497344.0
$ws.Range("I132").Value = 58830380.0
$ws.Range("K132").Value = 176491140.0
$ws.Range("M132").Value = -176488610.0
$ws.Range("H141").Value = 109499.8
$ws.Range("J141").Value = 109999.81
$ws.Range("L141").Value = 109999.81
$ws.Range("N141").Value = -120359.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 248.0
$ws.Range("J86").Value = 250.0
$ws.Range("L86").Value = 750.0
$ws.Range("N86").Value = -3122.0
$ws.Range("H89").Value = 248.0
$ws.Range("J89").Value = 250.0
$ws.Range("L89").Value = 2250.0
$ws.Range("N89").Value = -14106.0
$ws.Range("H120").Value = 19860.572
$ws.Range("I120").Value = 9515.0
$ws.Range("J120").Value = 23998.8
$ws.Range("K120").Value = 28545.0
$ws.Range("L120").Value = 71996.4
$ws.Range("M120").Value = -23707.0
$ws.Range("N120").Value = -81672.4
$ws.Range("H129").Value = 1743.0476
$ws.Range("I129").Value = 954.46155
$ws.Range("J129").Value = 3024.5
$ws.Range("K129").Value = 2863.38465
$ws.Range("L129").Value = 9073.5
$ws.Range("M129").Value = 2136.61535
$ws.Range("N129").Value = -19073.5
$ws.Range("H132").Value = 1135.2
$ws.Range("I132").Value = 1019.0769
$ws.Range("J132").Value = 1890.0
$ws.Range("K132").Value = 9171.6921
$ws.Range("L132").Value = 17010.0
$ws.Range("M132").Value = -6641.6921
$ws.Range("N132").Value = -22070.0
$ws.Range("H137").Value = 44738812.0
$ws.Range("J137").Value = 11113826.0
$ws.Range("L137").Value = 33341478.0
$ws.Range("N137").Value = -33351678.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2002320.0
$ws.Range("I11").Value = 805201.2
$ws.Range("J11").Value = 2354413.8
$ws.Range("K11").Value = 805201.2
$ws.Range("L11").Value = 2354413.8
$ws.Range("M11").Value = -805062.2
$ws.Range("N11").Value = -2354691.8
$ws.Range("H31").Value = 512.75
$ws.Range("I31").Value = 512.75
$ws.Range("K31").Value = 512.75
$ws.Range("M31").Value = -220.75
$ws.Range("H37").Value = 512.75
$ws.Range("I37").Value = 512.75
$ws.Range("K37").Value = 512.75
$ws.Range("M37").Value = -235.75
$ws.Range("H70").Value = 6962.1035
$ws.Range("I70").Value = 6888.9443
$ws.Range("J70").Value = 7081.8184
$ws.Range("K70").Value = 6888.9443
$ws.Range("L70").Value = 7081.8184
$ws.Range("M70").Value = -6618.9443
$ws.Range("N70").Value = -7621.8184
$ws.Range("H73").Value = 6962.1035
$ws.Range("I73").Value = 6888.9443
$ws.Range("J73").Value = 7081.8184
$ws.Range("K73").Value = 6888.9443
$ws.Range("L73").Value = 7081.8184
$ws.Range("M73").Value = -5952.9443
$ws.Range("N73").Value = -8953.8184
$ws.Range("H96").Value = 0.0
$ws.Range("J96").Value = 0.0
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 4727.375
$ws.Range("I113").Value = 4596.5
$ws.Range("K113").Value = 4596.5
$ws.Range("M113").Value = -2426.5
$ws.Range("H122").Value = 561498.75
$ws.Range("I122").Value = 1668963.9
$ws.Range("K122").Value = 5006891.699999999
$ws.Range("M122").Value = -5004441.699999999
$ws.Range("H123").Value = 32893.332
$ws.Range("J123").Value = 32893.332
$ws.Range("L123").Value = 32893.332
$ws.Range("N123").Value = -37793.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 15000.0
$ws.Range("J4").Value = 15000.0
$ws.Range("L4").Value = 15000.0
$ws.Range("N4").Value = -15226.0
$ws.Range("H28").Value = 15000.0
$ws.Range("J28").Value = 15000.0
$ws.Range("L28").Value = 15000.0
$ws.Range("N28").Value = -15464.0
$ws.Range("H37").Value = 15000.0
$ws.Range("J37").Value = 15000.0
$ws.Range("L37").Value = 15000.0
$ws.Range("N37").Value = -15214.0
$ws.Range("H100").Value = 1828.5714
$ws.Range("I100").Value = 1800.0
$ws.Range("K100").Value = 1800.0
$ws.Range("M100").Value = -1259.0
$ws.Range("H122").Value = 6745.9287
$ws.Range("I122").Value = 3127.9167
$ws.Range("J122").Value = 9459.4375
$ws.Range("K122").Value = 9383.750100000001
$ws.Range("L122").Value = 28378.3125
$ws.Range("M122").Value = -6933.750100000001
$ws.Range("N122").Value = -33278.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 0.0
$ws.Range("I100").Value = 0.0
$ws.Range("J100").Value = 0.0
$ws.Range("K100").Value = 0.0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 3168.4443
$ws.Range("I122").Value = 2022.697
$ws.Range("K122").Value = 6068.090999999999
$ws.Range("M122").Value = -3618.090999999999
$ws.Range("H138").Value = 111159.6
$ws.Range("J138").Value = 131449.5
$ws.Range("L138").Value = 131449.5
$ws.Range("N138").Value = -141729.5
